# Add the "Map" sheet evaluation data
# (Deliverables/Code Review/Evaluation.xlsx)
# Commit message: "added map eval, and created Documentation Folder"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map")

# Row 4 -> Code readability (only applicable to game)
$ws.Range("B4").Value = "N/A"
$ws.Range("C4").Value = "N/A"

# Row 5 -> Security in Word Press
$ws.Range("B5").Value = "Secure as the user input is handled by the plugin and the google api key is not accesable to the outside user ."
$ws.Range("C5").Value = "N/A"

# Row 6 -> Maintainability (column C entered before column B)
$ws.Range("C6").Value = "Make sure that the page is protected from bot spams or other attacks that could cause the API to be charged extra. "
$ws.Range("B6").Value = "Easy to maintan, more locations can be easily added, and as long as the number of hits to the page stay at a reasonable number there should be no charges for the API. The plugin can also be simply updated through the wp-admin page"

# Row 8 -> Not needlessly complex (entered before row 7)
$ws.Range("B8").Value = "Pretty simple setup, simply a plugin that manages location, and a google maps api key that provides the maps functionality "

# Row 7 -> Wordpress Plugins active support
$ws.Range("B7").Value = "You can submit support requests through https://wpstorelocator.co/support/"

# Row 9 -> Documentation
$ws.Range("B9").Value = "There Is doucmentation for how to manage the plugin and how it is connected to the google API https://wpstorelocator.co/documentation/ "
$ws.Range("C9").Value = "Documentation on the how the Google Maps API key is set up, how billing is handled and what account it is on should all be created. "

# Row 10 -> Successful Testing
$ws.Range("B10").Value = "Tested the map from mulitple devices and locations."
$ws.Range("C10").Value = "Test the map on multiple browsers "

# Row 11 -> Speed
$ws.Range("B11").Value = "Map works quickly "

# Row 14 -> Outstanding Bugs/Glitches (B before C, ahead of rows 12/13/15/16 content)
$ws.Range("B14").Value = "The map functions correctly"

# Row 15 -> Aesthetic Design
$ws.Range("B15").Value = "Currently a few white boxes and a map make up the interface."

# Row 16 -> Is responsive to multiple devices
$ws.Range("B16").Value = "Work well on the phone and desktop"

# Row 14 / 15 column C, entered after B14/B15/B16
$ws.Range("C14").Value = "There is a small annoyance, when you navigate to the page it asks your location and displays the gyms and parks near you, but to filter by gyms or parks you have to manually enter your location and search again. This will be worked on as a stretch Goal "
$ws.Range("C15").Value = "The interface could look a little better, maybe adding some round edges and making the filter options more clear. This will be worked on as a stretch Goal "

# Remaining cells reuse already-existing shared strings.
$ws.Range("C7").Value = "N/A"
$ws.Range("C8").Value = "N/A"
$ws.Range("C11").Value = "N/A"
$ws.Range("B12").Value = "Namecheap provides automatic backups in the cpanel."
$ws.Range("C12").Value = "N/A"
$ws.Range("B13").Value = "Contact the namecheap"
$ws.Range("C13").Value = "N/A"
$ws.Range("C16").Value = "N/A"

# Make the "Map" tab the active/selected one, matching the new tabSelected state,
# and move the sheet's own selection to C16.
$ws.Activate()
$ws.Range("C16").Select()
